$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2813447.2
$ws.Range("J17").Value = 2813447.2
$ws.Range("L17").Value = 8440341.600000001
$ws.Range("N17").Value = -8440677.600000001

$ws.Range("H57").Value = 38499.5
$ws.Range("J57").Value = 38499.5
$ws.Range("L57").Value = 115498.5
$ws.Range("N57").Value = -116496.5

$ws.Range("H129").Value = 2330.8333
$ws.Range("I129").Value = 1409.1428
$ws.Range("J129").Value = 2917.3635
$ws.Range("K129").Value = 4227.428400000001
$ws.Range("L129").Value = 8752.0905
$ws.Range("M129").Value = 772.5715999999993
$ws.Range("N129").Value = -18752.0905

$ws.Range("H133").Value = 135000
$ws.Range("J133").Value = 135000
$ws.Range("L133").Value = 135000
$ws.Range("N133").Value = -145120

$ws.Range("H137").Value = 4393.643
$ws.Range("I137").Value = 1524
$ws.Range("J137").Value = 17594
$ws.Range("K137").Value = 4572
$ws.Range("L137").Value = 52782
$ws.Range("M137").Value = -2022
$ws.Range("N137").Value = -57882

$ws.Range("H138").Value = 7085.523
$ws.Range("I138").Value = 1445.625
$ws.Range("J138").Value = 10308.321
$ws.Range("K138").Value = 4336.875
$ws.Range("L138").Value = 30924.963
$ws.Range("M138").Value = 803.125
$ws.Range("N138").Value = -41204.963

$ws.Range("H141").Value = 9383.615
$ws.Range("I141").Value = 12320.444
$ws.Range("J141").Value = 2775.75
$ws.Range("K141").Value = 36961.33199999999
$ws.Range("L141").Value = 8327.25
$ws.Range("M141").Value = -31781.33199999999
$ws.Range("N141").Value = -18687.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6337.6665
$ws.Range("J45").Value = 6337.6665
$ws.Range("L45").Value = 6337.6665
$ws.Range("N45").Value = -7091.6665

$ws.Range("H61").Value = 125002400
$ws.Range("I61").Value = 166668660
$ws.Range("K61").Value = 166668660
$ws.Range("M61").Value = -166668448

$ws.Range("H74").Value = 35717216
$ws.Range("I74").Value = 250001920
$ws.Range("J74").Value = 3099.9167
$ws.Range("K74").Value = 250001920
$ws.Range("L74").Value = 3099.9167
$ws.Range("M74").Value = -250001046
$ws.Range("N74").Value = -4847.9167

$ws.Range("H77").Value = 35717216
$ws.Range("I77").Value = 250001920
$ws.Range("J77").Value = 3099.9167
$ws.Range("K77").Value = 1250009600
$ws.Range("L77").Value = 15499.5835
$ws.Range("M77").Value = -1250005232
$ws.Range("N77").Value = -24235.5835

$ws.Range("H110").Value = 5025.72
$ws.Range("I110").Value = 5177.5713
$ws.Range("K110").Value = 5177.5713
$ws.Range("M110").Value = -3132.5713

$ws.Range("H122").Value = 20836382
$ws.Range("I122").Value = 2173.5
$ws.Range("J122").Value = 41670590
$ws.Range("K122").Value = 6520.5
$ws.Range("L122").Value = 125011770
$ws.Range("M122").Value = -4070.5
$ws.Range("N122").Value = -125016670

$ws.Range("H127").Value = 43993
$ws.Range("J127").Value = 43993
$ws.Range("L127").Value = 43993
$ws.Range("N127").Value = -53913

$ws.Range("H136").Value = 125002400
$ws.Range("I136").Value = 166668660
$ws.Range("K136").Value = 500005980
$ws.Range("M136").Value = -500003430

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 68520
$ws.Range("J50").Value = 68520
$ws.Range("L50").Value = 68520
$ws.Range("N50").Value = -69668

$ws.Range("H60").Value = 41314.5
$ws.Range("J60").Value = 41314.5
$ws.Range("L60").Value = 41314.5
$ws.Range("N60").Value = -42512.5

$ws.Range("H81").Value = 26693.75
$ws.Range("J81").Value = 26693.75
$ws.Range("L81").Value = 26693.75
$ws.Range("N81").Value = -28815.75

$ws.Range("H84").Value = 26693.75
$ws.Range("J84").Value = 26693.75
$ws.Range("L84").Value = 80081.25
$ws.Range("N84").Value = -90689.25

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H120").Value = 68000
$ws.Range("J120").Value = 68000
$ws.Range("L120").Value = 68000
$ws.Range("N120").Value = -77676

$ws.Range("H127").Value = 52390
$ws.Range("J127").Value = 52390
$ws.Range("L127").Value = 52390
$ws.Range("N127").Value = -62310

$ws.Range("H130").Value = 80562
$ws.Range("J130").Value = 80562
$ws.Range("L130").Value = 80562
$ws.Range("N130").Value = -90602

$ws.Range("H135").Value = 80172.5
$ws.Range("J135").Value = 80172.5
$ws.Range("L135").Value = 80172.5
$ws.Range("N135").Value = -90312.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9620237
$ws.Range("I31").Value = 3183.2144
$ws.Range("J31").Value = 20840132
$ws.Range("K31").Value = 3183.2144
$ws.Range("L31").Value = 20840132
$ws.Range("M31").Value = -2888.2144
$ws.Range("N31").Value = -20840722

$ws.Range("H34").Value = 9620237
$ws.Range("I34").Value = 3183.2144
$ws.Range("J34").Value = 20840132
$ws.Range("K34").Value = 3183.2144
$ws.Range("L34").Value = 20840132
$ws.Range("M34").Value = -2981.2144
$ws.Range("N34").Value = -20840536

$ws.Range("H53").Value = 127230.664
$ws.Range("I53").Value = 58000
$ws.Range("J53").Value = 135884.5
$ws.Range("K53").Value = 58000
$ws.Range("L53").Value = 135884.5
$ws.Range("M53").Value = -57393
$ws.Range("N53").Value = -137098.5

$ws.Range("H98").Value = 38900
$ws.Range("J98").Value = 38900
$ws.Range("L98").Value = 38900
$ws.Range("N98").Value = -43392

$ws.Range("H104").Value = 64285
$ws.Range("J104").Value = 64285
$ws.Range("L104").Value = 64285
$ws.Range("N104").Value = -69527

$ws.Range("H119").Value = 78000
$ws.Range("J119").Value = 78000
$ws.Range("L119").Value = 78000
$ws.Range("N119").Value = -87676

$ws.Range("H121").Value = 55500
$ws.Range("J121").Value = 37000
$ws.Range("L121").Value = 37000
$ws.Range("N121").Value = -39620

$ws.Range("H122").Value = 2107242.8
$ws.Range("I122").Value = 1944.3684
$ws.Range("K122").Value = 5833.1052
$ws.Range("M122").Value = -3383.1052

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H138").Value = 119312
$ws.Range("J138").Value = 65520
$ws.Range("L138").Value = 65520
$ws.Range("N138").Value = -75800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2000
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H131").Value = 1225.6154
$ws.Range("I131").Value = 806.6667
$ws.Range("J131").Value = 1796.909
$ws.Range("K131").Value = 2420.0001
$ws.Range("L131").Value = 5390.727000000001
$ws.Range("M131").Value = 2619.9999
$ws.Range("N131").Value = -15470.727

$ws.Range("H133").Value = 5967.75
$ws.Range("I133").Value = 1290.3334
$ws.Range("J133").Value = 20000
$ws.Range("K133").Value = 3871.0002
$ws.Range("L133").Value = 60000
$ws.Range("M133").Value = 1188.9998
$ws.Range("N133").Value = -70120

$ws.Range("H134").Value = 6424.75
$ws.Range("I134").Value = 1349.5
$ws.Range("J134").Value = 11500
$ws.Range("K134").Value = 4048.5
$ws.Range("L134").Value = 34500
$ws.Range("M134").Value = 1021.5
$ws.Range("N134").Value = -44640

$ws.Range("H136").Value = 2415
$ws.Range("I136").Value = 2415
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7245
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2145
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1041926.6
$ws.Range("J2").Value = 358.54544
$ws.Range("L2").Value = 358.54544
$ws.Range("N2").Value = -584.54544

$ws.Range("H108").Value = 76666
$ws.Range("J108").Value = 74998
$ws.Range("L108").Value = 74998
$ws.Range("N108").Value = -82678

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H128").Value = 55779.668
$ws.Range("J128").Value = 55779.668
$ws.Range("L128").Value = 55779.668
$ws.Range("N128").Value = -65739.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1790.1305
$ws.Range("J22").Value = 1748.6
$ws.Range("L22").Value = 1748.6
$ws.Range("N22").Value = -2338.6

$ws.Range("H27").Value = 1790.1305
$ws.Range("J27").Value = 1748.6
$ws.Range("L27").Value = 1748.6
$ws.Range("N27").Value = -1962.6

$ws.Range("H74").Value = 62000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 62000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H92").Value = 69990
$ws.Range("I92").Value = 69990
$ws.Range("K92").Value = 69990
$ws.Range("M92").Value = -67494

$ws.Range("H99").Value = 29129.5
$ws.Range("I99").Value = 29129.5
$ws.Range("K99").Value = 29129.5
$ws.Range("M99").Value = -26134.5

$ws.Range("H102").Value = 68749.75
$ws.Range("J102").Value = 73999.5
$ws.Range("L102").Value = 73999.5
$ws.Range("N102").Value = -80489.5

$ws.Range("H123").Value = 69565.42999999999
$ws.Range("J123").Value = 74993
$ws.Range("L123").Value = 74993
$ws.Range("N123").Value = -84793

$ws.Range("H129").Value = 60000
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H131").Value = 69836.75
$ws.Range("J131").Value = 87025.5
$ws.Range("L131").Value = 87025.5
$ws.Range("N131").Value = -97105.5

$ws.Range("H132").Value = 4685.143
$ws.Range("I132").Value = 2599.6667
$ws.Range("K132").Value = 7799.000100000001
$ws.Range("M132").Value = -5269.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 32997
$ws.Range("J70").Value = 32997
$ws.Range("L70").Value = 32997
$ws.Range("N70").Value = -33627

$ws.Range("H73").Value = 32997
$ws.Range("J73").Value = 32997
$ws.Range("L73").Value = 32997
$ws.Range("N73").Value = -35181

$ws.Range("H75").Value = 40496.668
$ws.Range("I75").Value = 34333.332
$ws.Range("J75").Value = 46660
$ws.Range("K75").Value = 34333.332
$ws.Range("L75").Value = 46660
$ws.Range("M75").Value = -33397.332
$ws.Range("N75").Value = -48532

$ws.Range("H78").Value = 40496.668
$ws.Range("I78").Value = 34333.332
$ws.Range("J78").Value = 46660
$ws.Range("K78").Value = 102999.996
$ws.Range("L78").Value = 139980
$ws.Range("M78").Value = -98319.99600000001
$ws.Range("N78").Value = -149340

$ws.Range("H86").Value = 12540624
$ws.Range("J86").Value = 40998.2
$ws.Range("L86").Value = 40998.2
$ws.Range("N86").Value = -43244.2

$ws.Range("H87").Value = 67500
$ws.Range("I87").Value = 67500
$ws.Range("K87").Value = 67500
$ws.Range("M87").Value = -66252

$ws.Range("H89").Value = 12540624
$ws.Range("J89").Value = 40998.2
$ws.Range("L89").Value = 204991
$ws.Range("N89").Value = -216223

$ws.Range("H90").Value = 67500
$ws.Range("I90").Value = 67500
$ws.Range("K90").Value = 202500
$ws.Range("M90").Value = -196260

$ws.Range("H102").Value = 69000
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H106").Value = 42000
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H127").Value = 60000
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H129").Value = 60112.5
$ws.Range("J129").Value = 70225
$ws.Range("L129").Value = 70225
$ws.Range("N129").Value = -80225

$ws.Range("H131").Value = 63357
$ws.Range("J131").Value = 63357
$ws.Range("L131").Value = 63357
$ws.Range("N131").Value = -73437

$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200
